$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date in column C for rows 2-5
# from 2023-09-01 (45170) to 2023-09-05 (45174)
$ws.Range("C2:C5").Value = (Get-Date -Year 2023 -Month 9 -Day 5).Date
